# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.
#
# Values that look like dates (e.g. "2026-02-06") or percentages
# (e.g. "67.7%") are entered with a leading apostrophe so Excel stores
# them as literal text (quote-prefixed), matching the plain-text values
# already used throughout these logs instead of converting them into
# real date/number values.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$row,
        [string]$date,
        [string]$timestamp,
        [string]$hour,
        [string]$location,
        [string]$value,
        [string]$status
    )

    $ws.Range("A$row").Value = "'" + $date
    $ws.Range("B$row").Value = $timestamp
    $ws.Range("C$row").Value = $hour
    $ws.Range("D$row").Value = $location

    # Percentage-looking values ("67.7%") must also be forced to text so
    # they don't get reinterpreted as numeric percentages.
    if ($value -match '^[0-9]+(\.[0-9]+)?%$') {
        $ws.Range("E$row").Value = "'" + $value
    } else {
        $ws.Range("E$row").Value = $value
    }

    $ws.Range("F$row").Value = $status
}

# ---------------------------------------------------------------------------
# PIR sheet: rows 607-614
# ---------------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")

Add-LogRow $pir 607 "2026-02-06" "10:30:23" "10:00" "Bathroom" "No Motion"       "Inactive"
Add-LogRow $pir 608 "2026-02-06" "10:30:25" "10:00" "Bathroom" "No Motion"       "Inactive"
Add-LogRow $pir 609 "2026-02-06" "10:30:29" "10:00" "Bathroom" "No Motion"       "Inactive"
Add-LogRow $pir 610 "2026-02-06" "10:30:33" "10:00" "Bathroom" "No Motion"       "Inactive"
Add-LogRow $pir 611 "2026-02-06" "10:30:39" "10:00" "Bathroom" "No Motion"       "Inactive"
Add-LogRow $pir 612 "2026-02-06" "10:30:43" "10:00" "Bathroom" "No Motion"       "Inactive"
Add-LogRow $pir 613 "2026-02-06" "10:30:49" "10:00" "Bathroom" "No Motion"       "Inactive"
Add-LogRow $pir 614 "2026-02-06" "10:30:50" "10:00" "Bathroom" "Motion Detected" "Active"

# ---------------------------------------------------------------------------
# Humidity sheet: rows 436-440
# ---------------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")

Add-LogRow $humidity 436 "2026-02-06" "10:30:26" "10:00" "Bathroom" "67.7%" "Active"
Add-LogRow $humidity 437 "2026-02-06" "10:30:30" "10:00" "Bathroom" "66.8%" "Active"
Add-LogRow $humidity 438 "2026-02-06" "10:30:40" "10:00" "Bathroom" "67.1%" "Active"
Add-LogRow $humidity 439 "2026-02-06" "10:30:45" "10:00" "Bathroom" "68.0%" "Active"
Add-LogRow $humidity 440 "2026-02-06" "10:30:54" "10:00" "Bathroom" "67.8%" "Active"

# ---------------------------------------------------------------------------
# Temperature sheet: rows 435-439
# ---------------------------------------------------------------------------
$temperature = $wb.Worksheets.Item("Temperature")

Add-LogRow $temperature 435 "2026-02-06" "10:30:22" "10:00" "Bathroom" "28.4C" "Active"
Add-LogRow $temperature 436 "2026-02-06" "10:30:27" "10:00" "Bathroom" "28.4C" "Active"
Add-LogRow $temperature 437 "2026-02-06" "10:30:31" "10:00" "Bathroom" "28.3C" "Active"
Add-LogRow $temperature 438 "2026-02-06" "10:30:41" "10:00" "Bathroom" "28.3C" "Active"
Add-LogRow $temperature 439 "2026-02-06" "10:30:46" "10:00" "Bathroom" "28.4C" "Active"
